# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme name="Office")
#   ppt/theme/theme2.xml -> "Integral"     (clrScheme name="Red Violet")
# theme2.xml is the theme actually wired to the slide master (and thus to
# every slide), so it's the one reachable from the PowerPoint object model
# via SlideMaster.Theme. This swaps the "Office" colour values into the
# live theme so the deck renders with the original Office palette instead
# of the Red Violet / Integral one.
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
$tcs = $t.ThemeColorScheme

# Colour order matches the clrScheme child order in the OOXML:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB integers are PowerPoint's 0xBBGGRR encoding of the target hex colour.
$tcs.Item(1).RGB  = 0         # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477   # folHlink -> 954F72
